$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Block 1: insert two blank paragraphs + an "Obs:" (red) paragraph about the
# Pareto Principle right after the paragraph that ends in
# "...Princípio de Pareto." (and before the "Cite 3 exemplos..." paragraph).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Dois exemplos são o uso de um Mapa Mental ou Princípio de Pareto.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara1 = $rng1.Paragraphs(1)

$block1 = "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:ind w:left='360'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='24292E'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:ind w:left='360'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
      "<w:t>Obs</w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
      "<w:t>: O Princípio de Pareto, ou regra 80/20, é uma tendência que prevê que 80% dos efeitos surgem a partir de apenas 20% das causas, podendo ser aplicado em várias outras relações de causa e efeito.</w:t>" +
    "</w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:ind w:left='360'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='24292E'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
  "</w:p>"

$anchorPara1.Range.InsertParagraphAfter()
$newPara1 = $anchorPara1.Next()
$newPara1.Range.InsertXML($block1)

# ---------------------------------------------------------------------------
# Block 2: insert a blank bullet paragraph plus three more paragraphs (two
# "Obs:" style red paragraphs, one blank) right after the paragraph that ends
# in "...Behavior Driven Development)" (and before the trailing blank
# paragraph / sectPr).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Behavior Driven Development)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara2 = $rng2.Paragraphs(1)

$block2 = "<w:p $wns>" +
    "<w:pPr>" +
      "<w:numPr>" +
        "<w:ilvl w:val='1'/>" +
        "<w:numId w:val='1'/>" +
      "</w:numPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='24292E'/>" +
        "<w:sz w:val='20'/>" +
        "<w:szCs w:val='20'/>" +
        "<w:lang w:eastAsia='pt-BR'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t>Obs</w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t>: A </w:t>" +
    "</w:r>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:b/>" +
        "<w:bCs/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t>modelagem</w:t>" +
    "</w:r>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t xml:space='preserve'> de processos é capaz de disponibilizar aos gestores a informação correta para a tomada de decisão no nível estratégico, gerenciar processos e garantir a coordenação das atividades no nível operacional.</w:t>" +
    "</w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:pPr>" +
      "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "<w:spacing w:after='0' w:line='360' w:lineRule='auto'/>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t xml:space='preserve'>O </w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:b/>" +
        "<w:bCs/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t>Gherkin</w:t>" +
    "</w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r>" +
      "<w:rPr>" +
        "<w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
        "<w:color w:val='FF0000'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
      "</w:rPr>" +
      "<w:t xml:space='preserve'> é um dos elementos principais quando se trata de BDD em automação. Sua função é padronizar a forma de descrever especificações de cenários, baseado na regra de negócio.</w:t>" +
    "</w:r>" +
  "</w:p>"

$anchorPara2.Range.InsertParagraphAfter()
$newPara2 = $anchorPara2.Next()
$newPara2.Range.InsertXML($block2)

Write-Host "Done"
